$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Change 1: append two new runs to the end of the first paragraph ---
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$rEnd1 = $d.Range($r1.End - 1, $r1.End - 1)
$rEnd1.InsertAfter(" It was the easiest to have work being done in")
$rEnd1.Collapse(0)
$rEnd1.InsertAfter(" every process at the same time, unlike blocking which would not always have working being done on every lock step.")

# --- Change 2: add five new list paragraphs after the final paragraph ---
# Remove the existing "_GoBack" bookmark up front; Word re-anchors this
# bookmark to the most-recently-edited location, which ends up being at the
# very end of the document once we've appended the new paragraphs below. We
# recreate the bookmark (same id/name) as part of inserting the last
# paragraph's XML.
if ($d.Bookmarks.Exists("_GoBack")) {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
}

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$insertPos = $lastPara.Range.End - 1

# Paragraph: "See Analytics excel sheet" (ilvl 0)
$xml1 = "<w:p $wNs><w:pPr><w:pStyle w:val=""ListParagraph""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""1""/></w:numPr></w:pPr><w:r><w:t>See Analytic</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space=""preserve""> excel sheet</w:t></w:r></w:p>"
$r = $d.Range($insertPos, $insertPos)
$r.InsertXML($xml1)

# Paragraph: " " (ilvl 0)
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$insertPos = $lastPara.Range.End - 1
$xml2 = "<w:p $wNs><w:pPr><w:pStyle w:val=""ListParagraph""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""1""/></w:numPr></w:pPr><w:r><w:t xml:space=""preserve""> </w:t></w:r></w:p>"
$r = $d.Range($insertPos, $insertPos)
$r.InsertXML($xml2)

# Paragraph: performance observations (ilvl 1)
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$insertPos = $lastPara.Range.End - 1
$xml3 = "<w:p $wNs><w:pPr><w:pStyle w:val=""ListParagraph""/><w:numPr><w:ilvl w:val=""1""/><w:numId w:val=""1""/></w:numPr></w:pPr>" +
        "<w:r><w:t xml:space=""preserve"">Performance increased from 1 to 2 processors, but decreased with 4 processors. My belief is that this has to do with how Mac OS was scheduling since I believe that I have 4 physical cores. As the matrix size increased, the performance increased. Having an additional </w:t></w:r>" +
        "<w:r><w:t xml:space=""preserve"">2 </w:t></w:r>" +
        "<w:r><w:t>processor</w:t></w:r>" +
        "<w:r><w:t>s</w:t></w:r>" +
        "<w:r><w:t xml:space=""preserve""> for </w:t></w:r>" +
        "<w:r><w:t>the large/ginormous matrix really seemed to improve things.</w:t></w:r>" +
        "</w:p>"
$r = $d.Range($insertPos, $insertPos)
$r.InsertXML($xml3)

# Paragraph: logarithmic decrease observation (ilvl 1)
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$insertPos = $lastPara.Range.End - 1
$xml4 = "<w:p $wNs><w:pPr><w:pStyle w:val=""ListParagraph""/><w:numPr><w:ilvl w:val=""1""/><w:numId w:val=""1""/></w:numPr></w:pPr><w:r><w:t>That the time decreased with 4 processors, I really thought it would be a logarithmic decrease.</w:t></w:r></w:p>"
$r = $d.Range($insertPos, $insertPos)
$r.InsertXML($xml4)

# Paragraph: MPI reflections (ilvl 1) -- also recreates the "_GoBack" bookmark
# at the very end of the document, matching its final position.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$insertPos = $lastPara.Range.End - 1
$xml5 = "<w:p $wNs><w:pPr><w:pStyle w:val=""ListParagraph""/><w:numPr><w:ilvl w:val=""1""/><w:numId w:val=""1""/></w:numPr></w:pPr><w:r><w:t>Doing the matrix multiply with MPI. I think had I implemented another MPI problem without matrix multiplication, I would have figured it out much faster because learning MPI took quite a bit of time. It`u2019s not very plug and play.</w:t></w:r><w:bookmarkStart w:id=""0"" w:name=""_GoBack""/><w:bookmarkEnd w:id=""0""/></w:p>"
$r = $d.Range($insertPos, $insertPos)
$r.InsertXML($xml5)
